# Driver code refactoring: the "password" column on the DATA sheet now
# stores a base64-encoded value instead of the previous plain-text one,
# and the old manual "no" test row for the Amazon hamburger-menu test
# case is removed (no longer needed after the refactor).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Encode the password values (admin123 -> YWRtaW4xMjM=) for every data
# row that previously held the plain-text password.
$ws.Range("E2:E6").Value = "YWRtaW4xMjM="

# Drop the last test row (row 9): amazonHamburgerMenuTest / no / chrome ...
$ws.Rows.Item(9).Delete()

# Keep the column reasonably sized for the longer encoded value.
$ws.Columns.Item(5).ColumnWidth = 14.91

# Move the active selection to where the last remaining row now is.
$null = $ws.Range("B7").Select()
